# Auto-generated market-data refresh for Sheets/Ultros_Profits.xlsx
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# per-row, per-sheet, mirroring the scheduled runner's data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 17.5
$ws.Range("I12").Value = 17.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 17.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 152.5
$ws.Range("N12").ClearContents()
$ws.Range("H28").Value = 847.2174
$ws.Range("I28").Value = 446.29413
$ws.Range("J28").Value = 1983.1666
$ws.Range("K28").Value = 446.29413
$ws.Range("L28").Value = 1983.1666
$ws.Range("M28").Value = 38.70587
$ws.Range("N28").Value = -2953.1666
$ws.Range("H33").Value = 162.92857
$ws.Range("I33").Value = 153.15384
$ws.Range("K33").Value = 153.15384
$ws.Range("M33").Value = 75.84616
$ws.Range("H40").Value = 3821.4517
$ws.Range("I40").Value = 3926.9167
$ws.Range("J40").Value = 3459.8572
$ws.Range("K40").Value = 3926.9167
$ws.Range("L40").Value = 3459.8572
$ws.Range("M40").Value = -3751.9167
$ws.Range("N40").Value = -3809.8572
$ws.Range("H41").Value = 938.2
$ws.Range("I41").Value = 771.1875
$ws.Range("J41").Value = 1606.25
$ws.Range("K41").Value = 771.1875
$ws.Range("L41").Value = 1606.25
$ws.Range("M41").Value = -331.1875
$ws.Range("N41").Value = -2486.25
$ws.Range("H43").Value = 4416.6665
$ws.Range("I43").Value = 4250
$ws.Range("K43").Value = 4250
$ws.Range("M43").Value = -4181
$ws.Range("H62").Value = 20166.334
$ws.Range("I62").Value = 28000
$ws.Range("K62").Value = 28000
$ws.Range("M62").Value = -27376
$ws.Range("H65").Value = 20166.334
$ws.Range("I65").Value = 28000
$ws.Range("K65").Value = 140000
$ws.Range("M65").Value = -136880
$ws.Range("H103").Value = 1139.25
$ws.Range("I103").Value = 1343.8572
$ws.Range("J103").Value = 852.8
$ws.Range("K103").Value = 4031.5716
$ws.Range("L103").Value = 2558.4
$ws.Range("M103").Value = -3445.5716
$ws.Range("N103").Value = -3730.4
$ws.Range("H113").Value = 10845.363
$ws.Range("I113").Value = 9566.5
$ws.Range("J113").Value = 12380
$ws.Range("K113").Value = 9566.5
$ws.Range("L113").Value = 12380
$ws.Range("M113").Value = -6312.5
$ws.Range("N113").Value = -18888
$ws.Range("H133").Value = 59769.23
$ws.Range("J133").Value = 59769.23
$ws.Range("L133").Value = 59769.23
$ws.Range("N133").Value = -69889.23000000001
$ws.Range("H138").Value = 4868.087
$ws.Range("I138").Value = 3499
$ws.Range("K138").Value = 10497
$ws.Range("M138").Value = -5357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13040.167
$ws.Range("I2").Value = 15339.583
$ws.Range("K2").Value = 15339.583
$ws.Range("M2").Value = -15226.583
$ws.Range("H97").Value = 3956.8333
$ws.Range("I97").Value = 2992.3157
$ws.Range("K97").Value = 2992.3157
$ws.Range("M97").Value = -2496.3157
$ws.Range("H102").Value = 1719.375
$ws.Range("I102").Value = 1567.3334
$ws.Range("K102").Value = 1567.3334
$ws.Range("M102").Value = 54.66660000000002
$ws.Range("H110").Value = 3380.6086
$ws.Range("I110").Value = 3226.0967
$ws.Range("K110").Value = 3226.0967
$ws.Range("M110").Value = -1181.0967
$ws.Range("H116").Value = 13040.167
$ws.Range("I116").Value = 15339.583
$ws.Range("K116").Value = 15339.583
$ws.Range("M116").Value = -13045.583
$ws.Range("H122").Value = 3754.2942
$ws.Range("I122").Value = 3201.6785
$ws.Range("J122").Value = 6333.1665
$ws.Range("K122").Value = 9605.0355
$ws.Range("L122").Value = 18999.4995
$ws.Range("M122").Value = -7155.0355
$ws.Range("N122").Value = -23899.4995
$ws.Range("H132").Value = 1498.05
$ws.Range("I132").Value = 1498.05
$ws.Range("K132").Value = 4494.15
$ws.Range("M132").Value = -1964.15

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13040.167
$ws.Range("I3").Value = 15339.583
$ws.Range("K3").Value = 15339.583
$ws.Range("M3").Value = -15225.583
$ws.Range("H9").Value = 44999.4
$ws.Range("I9").Value = 44999
$ws.Range("K9").Value = 44999
$ws.Range("M9").Value = -44831
$ws.Range("H20").Value = 113037.664
$ws.Range("I20").Value = 2048.4285
$ws.Range("J20").Value = 501500
$ws.Range("K20").Value = 2048.4285
$ws.Range("L20").Value = 501500
$ws.Range("M20").Value = -1801.4285
$ws.Range("N20").Value = -501994
$ws.Range("H94").Value = 2552.5625
$ws.Range("I94").Value = 2275.6924
$ws.Range("K94").Value = 2275.6924
$ws.Range("M94").Value = -1824.6924
$ws.Range("H105").Value = 4406
$ws.Range("I105").Value = 3594.2307
$ws.Range("K105").Value = 3594.2307
$ws.Range("M105").Value = -1847.2307
$ws.Range("H109").Value = 49857
$ws.Range("J109").Value = 53333.168
$ws.Range("L109").Value = 53333.168
$ws.Range("N109").Value = -56107.168
$ws.Range("H134").Value = 10446.947
$ws.Range("I134").Value = 2264.2354
$ws.Range("J134").Value = 80000
$ws.Range("K134").Value = 6792.706200000001
$ws.Range("L134").Value = 240000
$ws.Range("M134").Value = -4257.706200000001
$ws.Range("N134").Value = -245070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10179.615
$ws.Range("I31").Value = 15770.571
$ws.Range("J31").Value = 3656.8333
$ws.Range("K31").Value = 15770.571
$ws.Range("L31").Value = 3656.8333
$ws.Range("M31").Value = -15475.571
$ws.Range("N31").Value = -4246.8333
$ws.Range("H34").Value = 10179.615
$ws.Range("I34").Value = 15770.571
$ws.Range("J34").Value = 3656.8333
$ws.Range("K34").Value = 15770.571
$ws.Range("L34").Value = 3656.8333
$ws.Range("M34").Value = -15568.571
$ws.Range("N34").Value = -4060.8333
$ws.Range("H134").Value = 5330.1113
$ws.Range("I134").Value = 5369.625
$ws.Range("J134").Value = 5014
$ws.Range("K134").Value = 16108.875
$ws.Range("L134").Value = 15042
$ws.Range("M134").Value = -13573.875
$ws.Range("N134").Value = -20112

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 511.66666
$ws.Range("I12").Value = 355.66666
$ws.Range("J12").Value = 667.6667
$ws.Range("K12").Value = 1066.99998
$ws.Range("L12").Value = 2003.0001
$ws.Range("M12").Value = -893.9999800000001
$ws.Range("N12").Value = -2349.0001
$ws.Range("H14").Value = 166.55556
$ws.Range("I14").Value = 166.55556
$ws.Range("K14").Value = 499.66668
$ws.Range("M14").Value = -326.66668
$ws.Range("H68").Value = 1365.2559
$ws.Range("J68").Value = 1365.2559
$ws.Range("L68").Value = 4095.7677
$ws.Range("N68").Value = -5717.7677
$ws.Range("H71").Value = 1365.2559
$ws.Range("J71").Value = 1365.2559
$ws.Range("L71").Value = 12287.3031
$ws.Range("N71").Value = -20399.3031
$ws.Range("H98").Value = 1619.125
$ws.Range("I98").Value = 1263.5
$ws.Range("J98").Value = 1974.75
$ws.Range("K98").Value = 3790.5
$ws.Range("L98").Value = 5924.25
$ws.Range("M98").Value = -2292.5
$ws.Range("N98").Value = -8920.25
$ws.Range("H107").Value = 1812.8182
$ws.Range("I107").Value = 1299.9166
$ws.Range("J107").Value = 2428.3
$ws.Range("K107").Value = 3899.7498
$ws.Range("L107").Value = 7284.900000000001
$ws.Range("M107").Value = -1979.7498
$ws.Range("N107").Value = -11124.9
$ws.Range("H132").Value = 1038.2632
$ws.Range("I132").Value = 1044.4286
$ws.Range("J132").Value = 1021
$ws.Range("K132").Value = 9399.857399999999
$ws.Range("L132").Value = 9189
$ws.Range("M132").Value = -6869.857399999999
$ws.Range("N132").Value = -14249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 177.73334
$ws.Range("I2").Value = 60.909092
$ws.Range("K2").Value = 60.909092
$ws.Range("M2").Value = 52.090908
$ws.Range("H80").Value = 2439.1667
$ws.Range("I80").Value = 2439.1667
$ws.Range("K80").Value = 2439.1667
$ws.Range("M80").Value = -1441.1667
$ws.Range("H83").Value = 2439.1667
$ws.Range("I83").Value = 2439.1667
$ws.Range("K83").Value = 12195.8335
$ws.Range("M83").Value = -7203.833500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 13809.523
$ws.Range("J109").Value = 13809.523
$ws.Range("L109").Value = 13809.523
$ws.Range("N109").Value = -16583.523
$ws.Range("H132").Value = 4270.4688
$ws.Range("I132").Value = 4321.129
$ws.Range("K132").Value = 12963.387
$ws.Range("M132").Value = -10433.387

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 60777.223
$ws.Range("I96").Value = 170999
$ws.Range("J96").Value = 5666.3335
$ws.Range("K96").Value = 170999
$ws.Range("L96").Value = 5666.3335
$ws.Range("M96").Value = -169626
$ws.Range("N96").Value = -8412.333500000001
$ws.Range("H107").Value = 1071.1428
$ws.Range("I107").Value = 1100.8
$ws.Range("K107").Value = 3302.4
$ws.Range("M107").Value = -1382.4
$ws.Range("H122").Value = 4163.1665
$ws.Range("I122").Value = 3744.75
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 11234.25
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -8784.25
$ws.Range("N122").Value = -19900
$ws.Range("H125").Value = 49769.23
$ws.Range("J125").Value = 49769.23
$ws.Range("L125").Value = 49769.23
$ws.Range("N125").Value = -59609.23
$ws.Range("H136").Value = 8495.529
$ws.Range("I136").Value = 9833.333000000001
$ws.Range("K136").Value = 29499.999
$ws.Range("M136").Value = -26949.999
$ws.Range("H138").Value = 84666.664
$ws.Range("J138").Value = 84666.664
$ws.Range("L138").Value = 84666.664
$ws.Range("N138").Value = -94946.664

